# Atualização de bases das ligas, do dia: 16-05-2024 às 23:38
#
# The underlying match rows (id column B plus all stat columns C:AB) were
# re-sorted within a few same-date blocks. The sequential row id in column A
# stays put; everything from column B through AB moves between rows.
# Capture every source row's B:AB values BEFORE writing any of them back,
# since some rows feed each other (107<->108 swap, 142/143/144/145 is a
# 4-cycle, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(27, 28, 107, 108, 128, 129, 142, 143, 144, 145, 153, 154)

$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = $ws.Range("B$r`:AB$r").Value()
}

# target row -> source row (the data that now belongs in target used to live in source)
$mapping = @{
    27  = 28
    28  = 27
    107 = 108
    108 = 107
    128 = 129
    129 = 128
    142 = 143
    143 = 145
    144 = 142
    145 = 144
    153 = 154
    154 = 153
}

foreach ($target in $rows) {
    $source = $mapping[$target]
    $ws.Range("B$target`:AB$target").Value = $snapshot[$source]
}
